# feature: edit some payment
# Rename the "Credit Card" payment method to "CreditCard" (remove the space).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "CreditCard"
